$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2736504459652551
$ws.Range("C2").Value = 0.6164471910808867
$ws.Range("D2").Value = 0.5099410358193527
$ws.Range("E2").Value = 0.7141015584770507
$ws.Range("F2").Value = 0.6776623592093396
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = 0.1679866110975223
$ws.Range("C3").Value = 0.4776555476168131
$ws.Range("D3").Value = 0.3912275594536091
$ws.Range("E3").Value = 0.6254818618102439
$ws.Range("F3").Value = 0.6199689291560088
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.2407374606666166
$ws.Range("C4").Value = 0.3283263025081246
$ws.Range("D4").Value = 0.1549134947928615
$ws.Range("E4").Value = 0.393590516644979
$ws.Range("F4").Value = 0.320965582950402
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.3026896478973899
$ws.Range("C5").Value = 0.3981352914125425
$ws.Range("D5").Value = 0.2125958988159086
$ws.Range("E5").Value = 0.4610812280020827
$ws.Range("F5").Value = 0.3592211959175575
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.3313557959155047
$ws.Range("C6").Value = 0.3928030709233679
$ws.Range("D6").Value = 0.2145494597484939
$ws.Range("E6").Value = 0.463194839941567
$ws.Range("F6").Value = 0.3350151322728833
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.3718062059916627
$ws.Range("C7").Value = 0.3799713553340741
$ws.Range("D7").Value = 0.2141292450220557
$ws.Range("E7").Value = 0.4627410129025259
$ws.Range("F7").Value = 0.2858794074584022
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.4009368084426636
$ws.Range("C8").Value = 0.4127405556116039
$ws.Range("D8").Value = 0.227594880374374
$ws.Range("E8").Value = 0.4770690519981086
$ws.Range("F8").Value = 0.2691002335147636
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.4356636671697872
$ws.Range("C9").Value = 0.4368058703952601
$ws.Range("D9").Value = 0.2362224927271332
$ws.Range("E9").Value = 0.4860272551278716
$ws.Range("F9").Value = 0.2250325111912971
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.4167250561910649
$ws.Range("C10").Value = 0.4167250561910649
$ws.Range("D10").Value = 0.2156249438995315
$ws.Range("E10").Value = 0.4643543301182099
$ws.Range("F10").Value = 0.2148527137047464
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.3832484536769914
$ws.Range("C11").Value = 0.3832484536769914
$ws.Range("D11").Value = 0.1811246737723569
$ws.Range("E11").Value = 0.4255874455060403
$ws.Range("F11").Value = 0.1950649365568976
$ws.Range("G11").Value = 10
